$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 4166.6665
$ws.Range("J51").Value = 4166.6665
$ws.Range("L51").Value = 4166.6665
$ws.Range("N51").Value = -5134.6665
$ws.Range("H98").Value = 1168.5625
$ws.Range("I98").Value = 760.1
$ws.Range("K98").Value = 760.1
$ws.Range("M98").Value = 737.9
$ws.Range("H122").Value = 1168.5625
$ws.Range("I122").Value = 760.1
$ws.Range("K122").Value = 2280.3
$ws.Range("M122").Value = 169.6999999999998
$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").ClearContents()
$ws.Range("H141").Value = 4762.1113
$ws.Range("I141").Value = 4762.1113
$ws.Range("K141").Value = 14286.3339
$ws.Range("M141").Value = -9106.333899999998

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 224.90909
$ws.Range("I5").Value = 224.90909
$ws.Range("K5").Value = 224.90909
$ws.Range("M5").Value = -112.90909
$ws.Range("H32").Value = 4481.5356
$ws.Range("I32").Value = 4306.5
$ws.Range("K32").Value = 4306.5
$ws.Range("M32").Value = -4019.5
$ws.Range("H45").Value = 2058.8572
$ws.Range("I45").Value = 1985.3334
$ws.Range("K45").Value = 1985.3334
$ws.Range("M45").Value = -1608.3334
$ws.Range("H92").Value = 50000
$ws.Range("J92").Value = 50000
$ws.Range("L92").Value = 50000
$ws.Range("N92").Value = -54992
$ws.Range("H110").Value = 1934.4286
$ws.Range("I110").Value = 1799.2727
$ws.Range("K110").Value = 1799.2727
$ws.Range("M110").Value = 245.7273
$ws.Range("H122").Value = 1285
$ws.Range("I122").Value = 1285
$ws.Range("K122").Value = 3855
$ws.Range("M122").Value = -1405

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 224.90909
$ws.Range("I4").Value = 224.90909
$ws.Range("K4").Value = 224.90909
$ws.Range("M4").Value = -109.90909

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1198
$ws.Range("I16").Value = 1247.5
$ws.Range("K16").Value = 1247.5
$ws.Range("M16").Value = -960.5
$ws.Range("H31").Value = 2978.3103
$ws.Range("I31").Value = 2479.7058
$ws.Range("J31").Value = 3684.6667
$ws.Range("K31").Value = 2479.7058
$ws.Range("L31").Value = 3684.6667
$ws.Range("M31").Value = -2184.7058
$ws.Range("N31").Value = -4274.6667
$ws.Range("H34").Value = 2978.3103
$ws.Range("I34").Value = 2479.7058
$ws.Range("J34").Value = 3684.6667
$ws.Range("K34").Value = 2479.7058
$ws.Range("L34").Value = 3684.6667
$ws.Range("M34").Value = -2277.7058
$ws.Range("N34").Value = -4088.6667
$ws.Range("H58").Value = 1745.6364
$ws.Range("I58").Value = 1670.3
$ws.Range("K58").Value = 1670.3
$ws.Range("M58").Value = -1467.3
$ws.Range("H86").Value = 15000
$ws.Range("J86").Value = 15000
$ws.Range("L86").Value = 15000
$ws.Range("N86").Value = -17246
$ws.Range("H89").Value = 15000
$ws.Range("J89").Value = 15000
$ws.Range("L89").Value = 75000
$ws.Range("N89").Value = -86232
$ws.Range("H113").Value = 1198
$ws.Range("I113").Value = 1247.5
$ws.Range("K113").Value = 1247.5
$ws.Range("M113").Value = 922.5
$ws.Range("H136").Value = 1745.6364
$ws.Range("I136").Value = 1670.3
$ws.Range("K136").Value = 5010.9
$ws.Range("M136").Value = -2460.9

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 3000
$ws.Range("J17").Value = 3000
$ws.Range("L17").Value = 9000
$ws.Range("N17").Value = -9338
$ws.Range("H37").Value = 100000
$ws.Range("J37").Value = 100000
$ws.Range("L37").Value = 300000
$ws.Range("N37").Value = -300224
$ws.Range("H50").Value = 1639.8
$ws.Range("I50").Value = 999.6667
$ws.Range("K50").Value = 2999.0001
$ws.Range("M50").Value = -2518.0001
$ws.Range("H53").Value = 1639.8
$ws.Range("I53").Value = 999.6667
$ws.Range("K53").Value = 2999.0001
$ws.Range("M53").Value = -2518.0001

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 141.75
$ws.Range("I2").Value = 165.1
$ws.Range("J2").Value = 102.833336
$ws.Range("K2").Value = 165.1
$ws.Range("L2").Value = 102.833336
$ws.Range("M2").Value = -52.09999999999999
$ws.Range("N2").Value = -328.833336
$ws.Range("H122").Value = 3257.7693
$ws.Range("I122").Value = 2205
$ws.Range("K122").Value = 6615
$ws.Range("M122").Value = -4165
$ws.Range("H126").Value = 4176.4116
$ws.Range("I126").Value = 4066.6
$ws.Range("K126").Value = 12199.8
$ws.Range("M126").Value = -9729.799999999999
$ws.Range("H132").Value = 9931
$ws.Range("I132").Value = 9949.5
$ws.Range("J132").Value = 9894
$ws.Range("K132").Value = 29848.5
$ws.Range("L132").Value = 29682
$ws.Range("M132").Value = -27318.5
$ws.Range("N132").Value = -34742

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("M7").ClearContents()
$ws.Range("N7").ClearContents()
$ws.Range("H22").Value = 2256.8572
$ws.Range("I22").Value = 1799.6666
$ws.Range("J22").Value = 5000
$ws.Range("K22").Value = 1799.6666
$ws.Range("L22").Value = 5000
$ws.Range("M22").Value = -1504.6666
$ws.Range("N22").Value = -5590
$ws.Range("H27").Value = 2256.8572
$ws.Range("I27").Value = 1799.6666
$ws.Range("J27").Value = 5000
$ws.Range("K27").Value = 1799.6666
$ws.Range("L27").Value = 5000
$ws.Range("M27").Value = -1692.6666
$ws.Range("N27").Value = -5214
$ws.Range("H40").Value = 4666.6665
$ws.Range("I40").Value = 6000
$ws.Range("K40").Value = 6000
$ws.Range("M40").Value = -5864
$ws.Range("H105").Value = 30153.75
$ws.Range("J105").Value = 30153.75
$ws.Range("L105").Value = 30153.75
$ws.Range("N105").Value = -37141.75
$ws.Range("H122").Value = 3168.4443
$ws.Range("I122").Value = 3252.6667
$ws.Range("K122").Value = 9758.000100000001
$ws.Range("M122").Value = -7308.000100000001
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("N126").ClearContents()
$ws.Range("H132").Value = 4677.385
$ws.Range("I132").Value = 2472.2856
$ws.Range("J132").Value = 7250
$ws.Range("K132").Value = 7416.8568
$ws.Range("L132").Value = 21750
$ws.Range("M132").Value = -4886.8568
$ws.Range("N132").Value = -26810

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H21").Value = 312.5
$ws.Range("J21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("N21").ClearContents()
$ws.Range("H35").Value = 312.5
$ws.Range("J35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("N35").ClearContents()
$ws.Range("H100").Value = 1833
$ws.Range("I100").Value = 2499.5
$ws.Range("K100").Value = 4999
$ws.Range("M100").Value = -4458
$ws.Range("H126").Value = 1499
$ws.Range("I126").Value = 1499
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 4497
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -2027
$ws.Range("N126").ClearContents()
$ws.Range("H136").Value = 2106.2144
$ws.Range("I136").Value = 1581.091
$ws.Range("K136").Value = 4743.272999999999
$ws.Range("M136").Value = -2193.272999999999
